$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 0.225
$ws.Range("D3").Value = 0.575
$ws.Range("E3").Value = 0.9125
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = 0.366120218579235
$ws.Range("I3").Value = 0.1633304195804196
$ws.Range("J3").Value = 0.125
$ws.Range("K3").Value = 87.21250000000001
$ws.Range("T3").Value = 51
$ws.Range("U3").Value = 92
$ws.Range("Y3").Value = 235
$ws.Range("Z3").Value = 194
$ws.Range("AI3").Value = 0.821678
$ws.Range("AJ3").Value = 0.678322
